# Updates crypto price/symbol-list data per the "Updated symbol list" commit.
# Numeric-looking text values are written with a leading apostrophe so Excel
# keeps storing them as text (matching the workbook's existing inlineStr/text
# cell convention) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2 - BNB
Set-TextValue "D2" "244.82"

# Row 3 - OKB
Set-TextValue "D3" "23.81"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.321"

# Row 5 - Cronos
Set-TextValue "D5" "0.05777"

# Row 6
Set-TextValue "D6" "6.473"

# Row 8
Set-TextValue "D8" "0.8106"

# Row 9
Set-TextValue "D9" "0.8860"

# Row 10
Set-TextValue "D10" "0.1389"

# Row 11
Set-TextValue "D11" "0.07335"

# Row 12
Set-TextValue "D12" "0.03118"

# Row 13
Set-TextValue "D13" "0.03049"

# Row 14
Set-TextValue "D14" "0.09342"

# Row 15
Set-TextValue "D15" "3.861"

# Row 16
Set-TextValue "D16" "0.001561"

# Row 17
Set-TextValue "D17" "0.04722"

# Row 18
Set-TextValue "D18" "0.0006023"

# Row 19
Set-TextValue "D19" "0.005864"

# Row 20
Set-TextValue "D20" "0.001289"

# Row 22 - NitroEx
Set-TextValue "D22" "0.00008807"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"

# Row 23 - LEO
Set-TextValue "D23" "3.603"

# Row 26
Set-TextValue "D26" "0.1318"

# Row 41 - KickToken
Set-TextValue "D41" "0.006405"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 and 43 - CEJI / BKEXToken swapped identities, independent price updates
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1057"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002752"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008354"

# Row 45
Set-TextValue "D45" "0.00005391"

# Row 46
Set-TextValue "D46" "0.00000000750"

# Row 47
Set-TextValue "D47" "0.6903"

# Row 48
Set-TextValue "D48" "0.001844"

# Row 49
Set-TextValue "D49" "0.00002101"

Write-Output "done"
